$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, move the ResourceTag / ResourceType enum header block out of the
# way (rows 10 & 13 in column B) down to rows 48 & 51, since rows 10-13
# are being reused for the growing EffectType value list in columns G/J.
$ws.Range("B10").Copy()
$ws.Range("B48").PasteSpecial(-4122)
$ws.Range("B48").Value = "ResourceTag"
$ws.Range("B10").Clear()

$ws.Range("B13").Copy()
$ws.Range("B51").PasteSpecial(-4122)
$ws.Range("B51").Value = "ResourceType"
$ws.Range("B13").Clear()

# Seed the formatting for the new G8:G13 / J6:J13 cells by copying the
# existing formatted cell (G7) over the target range, then overwrite the
# values below (in the same order they were originally typed, so the
# shared-string table comes out in the same sequence).
$ws.Range("G7").Copy()
$ws.Range("G8:G13").PasteSpecial(-4122)
$ws.Range("J6:J13").PasteSpecial(-4122)

# New EffectType enum values (column G) for the DIE / DIESON rows.
$ws.Range("G8").Value  = "DIE"
$ws.Range("G9").Value  = "DIESON"

# Chinese comment/description column (column J) for each EffectType value.
$ws.Range("J6").Value  = "属性成长"
$ws.Range("J7").Value  = "结婚"
$ws.Range("J8").Value  = "死亡"
$ws.Range("J9").Value  = "孩子死亡"
$ws.Range("J10").Value = "获得资源"
$ws.Range("J11").Value = "获得钱"

# Remaining new EffectType enum values (column G).
$ws.Range("G12").Value = "ADD_CURR_EVENT"
$ws.Range("G11").Value = "ADD_MONEY"
$ws.Range("G10").Value = "ADD_RESOURCE"
$ws.Range("G13").Value = "ADD_NEXT_EVENT"

# Final two comments.
$ws.Range("J12").Value = "这一代获得新的事件"
$ws.Range("J13").Value = "下一代获得新的事件"

# GROWTH / MARRY already existed at G6 / G7 with the correct shared-string
# values, so they are left untouched (re-assigning is harmless too, since
# the runtime reuses the existing shared-string index for identical text).

# Match the final selection recorded in the saved workbook.
$ws.Range("H19").Select() | Out-Null
